$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.988.63"
$ws.Range("E2").Value = "  +5.71%  "
$ws.Range("D3").Value = "3.640.15"
$ws.Range("E3").Value = "  +16.31%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.59"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.98"
$ws.Range("E6").Value = "  +2.56%  "
$ws.Range("D7").Value = "3.637.25"
$ws.Range("E7").Value = "  +16.20%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  +3.41%  "
$ws.Range("E10").Value = "  +6.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.61"
$ws.Range("E11").Value = "  +3.01%  "
$ws.Range("E12").Value = "  +4.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.61"
$ws.Range("E13").Value = "  +11.52%  "
$ws.Range("E14").Value = "  +4.32%  "
$ws.Range("D15").Value = "4.247.68"
$ws.Range("E15").Value = "  +16.36%  "
$ws.Range("D16").Value = "70.932.90"
$ws.Range("E16").Value = "  +5.78%  "
$ws.Range("D17").Value = "3.647.35"
$ws.Range("E17").Value = "  +16.46%  "
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.47"
$ws.Range("E19").Value = "  +5.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.95"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "513.09"
$ws.Range("E21").Value = "  +4.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.14"
$ws.Range("E22").Value = "  +16.47%  "
$ws.Range("E23").Value = "  +6.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.31"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.48"
$ws.Range("E25").Value = "  +9.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.47"
$ws.Range("E26").Value = "  +4.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.98"
$ws.Range("E27").Value = "  +6.77%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.52"
$ws.Range("E29").Value = "  +9.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.17"
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("E31").Value = "  +17.00%  "
$ws.Range("E32").Value = "  +6.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.42"
$ws.Range("E33").Value = "  +11.38%  "
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  +7.71%  "
$ws.Range("E37").Value = "  +6.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.345"
$ws.Range("E38").Value = "  +10.76%  "
$ws.Range("E39").Value = "  +7.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.92"
$ws.Range("E40").Value = "  +2.93%  "
$ws.Range("E41").Value = "  +4.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.18"
$ws.Range("E42").Value = "  -7.13%  "
$ws.Range("D43").Value = "3.122.30"
$ws.Range("E43").Value = "  +11.00%  "
$ws.Range("E44").Value = "  +5.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "415.05"
$ws.Range("E45").Value = "  +10.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.80"
$ws.Range("E46").Value = "  +4.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.42"
$ws.Range("E47").Value = "  +14.05%  "
$ws.Range("E48").Value = "  +5.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "137.75"
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("E51").Value = "  +10.15%  "
